## Generate Report for Handoff
## Updates the localization-status workbook with the freshly generated
## handoff artifacts: the source doc GUID changed from 76aebd2e... to
## ea948818..., new xliff hashes were produced, and the handoff
## timestamps advanced.

$wb = $excel.ActiveWorkbook

$oldGuid = "76aebd2e-eb2d-4e6c-99dd-01a2c75bc336"
$newGuid = "ea948818-6ad9-446f-b6b8-e4bae19996e3"

$oldZhHash = "ae577d12b22dbeda2c7e4e83ee38e147e904de7b"
$newZhHash = "19db4c60dc66737a9cde3d0b15c459e1c8b97077"

$oldDeHash = "ae577d12b22dbeda2c7e4e83ee38e147e904de7b"
$newDeHash = "19db4c60dc66737a9cde3d0b15c459e1c8b97077"

$newMdName   = "$newGuid.md"
$newMdPath   = "e2e\$newGuid.md"
$newZhXlf    = "$newGuid.$newZhHash.zh-cn.xlf"
$newDeXlf    = "$newGuid.$newDeHash.de-de.xlf"

$newHoGenerateDate = "2016-08-27 00:57:32"
$newZhHandoffDate  = "2016-08-27 00:57:27"

## The hyperlink target (a GitHub blob URL pinned to a commit SHA) is
## unchanged by this edit - only the cell text / displayed caption moves
## to the new GUID. The host doesn't surface the existing Address via
## Hyperlinks.Item(N).Address (reads back empty), so it's restated here
## from the workbook's existing relationship.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f99f9c9fe4b369230f2d498f54844a9115e1b36/e2e/$oldGuid.md"
$linkColor = 0xED9564  ## RGB(0x64,0x95,0xED) - matches the workbook's existing HyperLink style

## ---- Overview sheet ----
$wsOverview = $wb.Sheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName

$linkRange = $wsOverview.Range("B2")
$linkRange.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($linkRange, $linkAddress, "", "", $newMdPath)
$linkRange.Font.Underline = $true
$linkRange.Font.Color = $linkColor

$wsOverview.Range("G2").Value = $newHoGenerateDate

## ---- zh-cn sheet ----
$wsZh = $wb.Sheets.Item("zh-cn")

$zhLinkRange = $wsZh.Range("A2")
$zhLinkRange.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($zhLinkRange, $linkAddress, "", "", $newMdName)
$zhLinkRange.Font.Underline = $true
$zhLinkRange.Font.Color = $linkColor

$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate

## ---- de-de sheet ----
$wsDe = $wb.Sheets.Item("de-de")

$deLinkRange = $wsDe.Range("A2")
$deLinkRange.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($deLinkRange, $linkAddress, "", "", $newMdName)
$deLinkRange.Font.Underline = $true
$deLinkRange.Font.Color = $linkColor

$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHoGenerateDate
